# Auto-update script: append the latest day of data to each sheet.
# - "Prix Spot": add a new date column BC ("07-aug") with 24 hourly prices.
# - "Gaz": add a new row (2025-08-05) with its closing price.
# - "CO2": add a new row (2025-08-05); price not published yet (blank).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": new column BC
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell, copy the format of the previous day's header (BB1) so the
# new column keeps the same bold/centered/bordered style, then set the text.
$wsPrix.Range("BB1").Copy()
$wsPrix.Range("BC1").PasteSpecial(-4122)
$wsPrix.Range("BC1").Value = "07-aug"

# Hourly prices for 07-aug, in row order (row 2 = "00 - 01" ... row 25 = "23 - 24").
$pricesBC = @(
    88.42,
    61.2,
    44.63,
    52.66,
    46.6,
    41.73,
    42.97,
    71.97,
    77.02,
    71.65000000000001,
    22.23,
    0,
    -0.01,
    -0.02,
    -0.07000000000000001,
    -0.01,
    3.81,
    48.97,
    70,
    101.53,
    115.42,
    110.74,
    110,
    101.28
)

for ($i = 0; $i -lt $pricesBC.Count; $i++) {
    $wsPrix.Cells.Item($i + 2, 55).Value = $pricesBC[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": new row 52
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Keep the date as plain text (matches every other row) instead of letting
# it be auto-converted to a date serial number.
$wsGaz.Range("A52").NumberFormat = "@"
$wsGaz.Range("A52").Value = "2025-08-05"
$wsGaz.Range("A51").Copy()
$wsGaz.Range("A52").PasteSpecial(-4122)

$wsGaz.Range("B52").Value = 33.775

# ---------------------------------------------------------------------------
# Sheet "CO2": new row 52 (price not available yet -> blank)
# ---------------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A52").NumberFormat = "@"
$wsCO2.Range("A52").Value = "2025-08-05"
$wsCO2.Range("A51").Copy()
$wsCO2.Range("A52").PasteSpecial(-4122)
